$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workload")

# Update the student nrs in row 5
$ws.Range("D5").Value = 4829360
$ws.Range("E5").Value = 5079934
$ws.Range("F5").Value = 5336724

# Update the selected cell/range on the sheet
$ws.Range("F8").Select()
